$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.397.37'
$ws.Range("E2").Value = '  -3.12%  '

$ws.Range("D3").Value = '3.375.61'
$ws.Range("E3").Value = '  -3.85%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.28'
$ws.Range("E5").Value = '  -3.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.09'
$ws.Range("E6").Value = '  -7.11%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '3.375.75'
$ws.Range("E8").Value = '  -3.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").Value = '  -2.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.24'
$ws.Range("E10").Value = '  -5.11%  '

$ws.Range("E11").Value = '  -4.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -3.58%  '

$ws.Range("D13").Value = '3.945.18'
$ws.Range("E13").Value = '  -4.04%  '

$ws.Range("E14").Value = '  -1.02%  '

$ws.Range("D15").Value = '3.366.92'
$ws.Range("E15").Value = '  -4.17%  '

$ws.Range("E16").Value = '  -6.21%  '

$ws.Range("D17").Value = '62.404.73'
$ws.Range("E17").Value = '  -3.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.31'
$ws.Range("E18").Value = '  -5.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.23'
$ws.Range("E19").Value = '  -7.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.62'
$ws.Range("E20").Value = '  -2.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.05'
$ws.Range("E21").Value = '  -4.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '372.37'
$ws.Range("E22").Value = '  -5.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("E23").Value = '  -4.59%  '

$ws.Range("D24").Value = '3.508.51'
$ws.Range("E24").Value = '  -3.89%  '

$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.46'
$ws.Range("E26").Value = '  -4.27%  '

$ws.Range("E27").Value = '  -10.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.89'
$ws.Range("E29").Value = '  -7.11%  '

$ws.Range("E30").Value = '  -7.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.76'
$ws.Range("E31").Value = '  -6.07%  '

$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.37'
$ws.Range("E33").Value = '  -7.05%  '

$ws.Range("D34").Value = '3.401.13'
$ws.Range("E34").Value = '  -3.90%  '

$ws.Range("E35").Value = '  -6.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.62'
$ws.Range("E36").Value = '  -3.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.17'
$ws.Range("E37").Value = '  -3.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '166.37'
$ws.Range("E38").Value = '  -0.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.61'
$ws.Range("E39").Value = '  -5.18%  '

$ws.Range("E40").Value = '  -5.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0749'
$ws.Range("E41").Value = '  -4.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.763'
$ws.Range("E43").Value = '  -5.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.47'
$ws.Range("E44").Value = '  -1.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.22'
$ws.Range("E45").Value = '  -5.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.51'
$ws.Range("E46").Value = '  -10.31%  '

$ws.Range("E47").Value = '  -7.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  -9.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.57'
$ws.Range("E49").Value = '  -3.43%  '

$ws.Range("D50").Value = '2.238.78'
$ws.Range("E50").Value = '  -5.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.835'
$ws.Range("E51").Value = '  -8.38%  '
